$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 1542
$ws.Range("L2").Value = 465
$ws.Range("M2").Value = 32
$ws.Range("N2").Value = 274
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 170
$ws.Range("T2").Value = 287
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 2471
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 2532
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 27
$ws.Range("AA2").Value = 15
